# Auto-generated edit script: updates cryptos list values (price/volume columns)
# to match the Sun Oct 22 06:55:35 UTC 2023 GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.110.98"

# Row 3
$ws.Range("D3").Value = "1.644.37"
$ws.Range("E3").Value = "  +2.65%  "

# Row 4
$ws.Range("E4").Value = "  +0.24%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.87"
$ws.Range("E5").Value = "  +1.65%  "

# Row 6
$ws.Range("E6").Value = "  +1.25%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.19%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "29.35"
$ws.Range("E8").Value = "  +4.75%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.261"
$ws.Range("E9").Value = "  +3.77%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0615"
$ws.Range("E10").Value = "  +2.04%  "

# Row 11
$ws.Range("E11").Value = "  +0.43%  "

# Row 12
$ws.Range("D12").Value = "1.879.72"
$ws.Range("E12").Value = "  +2.69%  "

# Row 13
$ws.Range("D13").Value = "1.648.14"
$ws.Range("E13").Value = "  +3.02%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.575"
$ws.Range("E14").Value = "  +5.41%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "9.47"
$ws.Range("E15").Value = "  +22.72%  "

# Row 16
$ws.Range("E16").Value = "  +4.63%  "

# Row 17
$ws.Range("D17").Value = "30.132.89"
$ws.Range("E17").Value = "  +1.65%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.96"
$ws.Range("E18").Value = "  +1.76%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "248.05"
$ws.Range("E19").Value = "  +2.37%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0711"
$ws.Range("E20").Value = "  +2.12%  "

# Row 21
$ws.Range("E21").Value = "  +0.13%  "

# Row 22
$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.01"
$ws.Range("E22").Value = "  +6.86%  "

# Row 23
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.21"
$ws.Range("E23").Value = "  +5.00%  "

# Row 24
$ws.Range("E24").Value = "  +2.07%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.22"
$ws.Range("E25").Value = "  +2.77%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.80"
$ws.Range("E26").Value = "  +2.30%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.111"
$ws.Range("E27").Value = "  +2.76%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.68"
$ws.Range("E28").Value = "  +3.94%  "

# Row 29
$ws.Range("E29").Value = "  +0.19%  "

# Row 30
$ws.Range("E30").Value = "  +2.71%  "

# Row 31
$ws.Range("E31").Value = "  +6.48%  "

# Row 32
$ws.Range("E32").Value = "  +6.17%  "

# Row 33
$ws.Range("E33").Value = "  +1.06%  "

# Row 34
$ws.Range("E34").Value = "  +1.18%  "

# Row 35
$ws.Range("E35").Value = "  +7.27%  "

# Row 36
$ws.Range("E36").Value = "  +1.98%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.87"
$ws.Range("E37").Value = "  -0.83%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "79.00"
$ws.Range("E38").Value = "  +19.44%  "

# Row 39
$ws.Range("E39").Value = "  +2.11%  "

# Row 40
$ws.Range("E40").Value = "  +0.12%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.560"
$ws.Range("E41").Value = "  +2.71%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.03"
$ws.Range("E42").Value = "  +2.86%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.849"
$ws.Range("E43").Value = "  +3.86%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "55.64"
$ws.Range("E44").Value = "  -2.50%  "

# Row 45
$ws.Range("E45").Value = "  +0.87%  "

# Row 46
$ws.Range("E46").Value = "  +5.65%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("E47").Value = "  +0.17%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.41"
$ws.Range("E48").Value = "  +1.18%  "

# Row 49
$ws.Range("D49").Value = "1.786.08"
$ws.Range("E49").Value = "  +2.60%  "

# Row 50
$ws.Range("E50").Value = "  +11.48%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "90.57"
$ws.Range("E51").Value = "  +4.38%  "

